$d = $word.ActiveDocument

# The document currently ends with a paragraph containing an inline image
# and the _GoBack bookmark. We need to append three new paragraphs after
# it: one empty paragraph, then two paragraphs of text.

$end = $d.Content.End

# Collapse to the very end of the document body (after the last paragraph
# mark) and insert the new paragraphs there.
$r = $d.Range($end, $end)

$r.InsertParagraphAfter()
$r.Collapse(0)

$r.InsertAfter("Le he pesto el AT holderName a la creditCard")
$r.InsertParagraphAfter()
$r.Collapse(0)

$r.InsertAfter("Le he puesto el max(100) al vat")
